# Insert a new row 37 (shifts existing rows 37-40 down to 38-41)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(37).Insert()

# Populate the new row 37 with the weekly entry (week of 2023-06-29, serial 45106)
$ws.Range("A37").Value = 10
$ws.Range("B37").Value = "Vega Modelo de Temuco"
$ws.Range("C37").Value = "La Araucanía"
$ws.Range("D37").Value = 45106
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100108
$ws.Range("H37").Value = "Tropicales y subtropicales"
$ws.Range("I37").Value = 100108001
$ws.Range("J37").Value = "Guayaba"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 120
$ws.Range("N37").Value = 2600
$ws.Range("O37").Value = 2600
$ws.Range("P37").Value = 2600
$ws.Range("Q37").Value = '$/kilo'
$ws.Range("R37").Value = "Región de Arica y Parinacota"
$ws.Range("S37").Value = 2600
$ws.Range("T37").Value = 1
